# Insert a new data row at row 448 (pushing the existing rows 448-512 down to
# 449-513) and populate it with the new Repollo price-sheet entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 448; this shifts rows 448:512
# down to 449:513 and extends the sheet dimension automatically.
$ws.Rows("448:448").Insert()

# Copy the date cell's number format from the row below (now 449) so the new
# date value renders the same way as the rest of column D.
$ws.Range("D448").NumberFormat = $ws.Range("D449").NumberFormat

$ws.Range("A448").Value2 = 5
$ws.Range("B448").Value2 = "Macroferia Regional de Talca"
$ws.Range("C448").Value2 = "Maule"
$ws.Range("D448").Value2 = 45127
$ws.Range("E448").Value2 = 7
$ws.Range("F448").Value2 = 100112006
$ws.Range("G448").Value2 = "Repollo"
$ws.Range("H448").Value2 = "Crespo record"
$ws.Range("I448").Value2 = "Primera"
$ws.Range("J448").Value2 = 5000
$ws.Range("K448").Value2 = 600
$ws.Range("L448").Value2 = 600
$ws.Range("M448").Value2 = 600
$ws.Range("N448").Value2 = "`$/unidad"
$ws.Range("O448").Value2 = "Región del Maule"
$ws.Range("P448").Value2 = 600
$ws.Range("Q448").Value2 = 1
$ws.Range("R448").Value2 = "Hortaliza"
